$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 295 (pushing old rows 295-394 down to 296-395),
# matching the weekly-refresh update described by the commit message.
$ws.Rows("295").Insert()

$ws.Range("A295").Value2 = 9
$ws.Range("B295").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C295").Value2 = "Metropolitana"
$ws.Range("D295").Value2 = 44988
$ws.Range("E295").Value2 = 13
$ws.Range("F295").Value2 = 300000001
$ws.Range("G295").Value2 = "Rabanito"
$ws.Range("H295").Value2 = "Sin especificar"
$ws.Range("I295").Value2 = "Primera"
$ws.Range("J295").Value2 = 7000
$ws.Range("K295").Value2 = 3000
$ws.Range("L295").Value2 = 3000
$ws.Range("M295").Value2 = 3000
$ws.Range("N295").Value2 = "$/cien unidades (volumen en unidades)"
$ws.Range("O295").Value2 = "Provincia de Chacabuco"
$ws.Range("P295").Value2 = 30
$ws.Range("Q295").Value2 = 100
$ws.Range("R295").Value2 = "Hortaliza"
